# Adds two new columns, I ("I0") and J ("IF"), to the right of the
# existing "IP" column (H), filling in the header labels and the
# per-row numeric values for rows 2-29.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Headers (row 1) -------------------------------------------------
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the header formatting (bold font, borders, centered alignment)
# from the existing "IP" header cell (H1) onto the two new header
# cells so they match the rest of the header row.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# --- Data rows (2-29): columns I (I0) and J (IF) ---------------------
$values = @{
    2  = @(1,6)
    3  = @(1,4)
    4  = @(1,5)
    5  = @(1,5)
    6  = @(1,6)
    7  = @(1,6)
    8  = @(1,7)
    9  = @(1,5)
    10 = @(1,4)
    11 = @(1,5)
    12 = @(1,5)
    13 = @(1,5)
    14 = @(8,8)
    15 = @(8,8)
    16 = @(7,8)
    17 = @(8,8)
    18 = @(7,8)
    19 = @(5,7)
    20 = @(9,9)
    21 = @(1,3)
    22 = @(1,7)
    23 = @(1,6)
    24 = @(1,3)
    25 = @(1,6)
    26 = @(1,5)
    27 = @(1,5)
    28 = @(1,3)
    29 = @(1,4)
}

foreach ($row in $values.Keys) {
    $pair = $values[$row]
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
}
